$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "29.942.42";    E = "  +0.10%  " },
    @{ Row = 3;  D = "1.895.41";     E = "  -0.15%  " },
    @{ Row = 4;  D = $null;          E = "  -0.04%  " },
    @{ Row = 5;  D = "0.7778";       E = "  -2.04%  " },
    @{ Row = 6;  D = "244.74";       E = "  +0.22%  " },
    @{ Row = 8;  D = "0.3137";       E = "  -1.07%  " },
    @{ Row = 9;  D = "25.90";        E = "  +1.53%  " },
    @{ Row = 10; D = "0.07262";      E = "  +1.07%  " },
    @{ Row = 11; D = "0.09300";      E = "  +14.63%  " },
    @{ Row = 12; D = "0.7747";       E = "  +0.55%  " },
    @{ Row = 13; D = "5.465";        E = "  -3.70%  " },
    @{ Row = 14; D = "1.905.24";     E = "  +0.63%  " },
    @{ Row = 15; D = "94.83";        E = "  +2.26%  " },
    @{ Row = 16; D = "6.233";        E = "  +0.79%  " },
    @{ Row = 17; D = "29.943.15";    E = "  +0.07%  " },
    @{ Row = 18; D = $null;          E = "  +0.06%  " },
    @{ Row = 19; D = "247.09";       E = "  +0.75%  " },
    @{ Row = 20; D = "0.000007892";  E = "  +1.54%  " },
    @{ Row = 21; D = "2.169.17";     E = "  +0.58%  " },
    @{ Row = 22; D = "8.172";        E = "  -0.81%  " },
    @{ Row = 23; D = $null;          E = "  -0.14%  " },
    @{ Row = 24; D = $null;          E = "  -0.06%  " },
    @{ Row = 25; D = "0.1592";       E = "  -4.85%  " },
    @{ Row = 26; D = "9.558";        E = "  +0.57%  " },
    @{ Row = 27; D = "162.55";       E = "  -1.08%  " },
    @{ Row = 28; D = "18.84";        E = "  +0.42%  " },
    @{ Row = 29; D = "2.052";        E = "  -1.11%  " },
    @{ Row = 30; D = "1.424";        E = "  +1.27%  " },
    @{ Row = 31; D = $null;          E = "  +0.17%  " },
    @{ Row = 32; D = "4.551";        E = "  +0.99%  " },
    @{ Row = 33; D = "4.126";        E = "  +0.83%  " },
    @{ Row = 34; D = "0.05526";      E = "  -1.84%  " },
    @{ Row = 35; D = $null;          E = "  -2.89%  " },
    @{ Row = 36; D = "0.7554";       E = "  +1.45%  " },
    @{ Row = 37; D = "1.001";        E = "  +0.01%  " },
    @{ Row = 38; D = "2.708";        E = "  +2.83%  " },
    @{ Row = 39; D = "0.01969";      E = "  +1.59%  " },
    @{ Row = 41; D = "0.4509";       E = $null },
    @{ Row = 42; D = "74.29";        E = "  -1.03%  " },
    @{ Row = 43; D = "6.096";        E = "  +2.44%  " },
    @{ Row = 44; D = "1.089.17";     E = "  -6.27%  " },
    @{ Row = 45; D = "0.8554";       E = "  +0.19%  " },
    @{ Row = 46; D = "1.000";        E = "  -0.04%  " },
    @{ Row = 47; D = "1.896";        E = "  +0.45%  " },
    @{ Row = 48; D = "102.80";       E = "  -1.71%  " },
    @{ Row = 49; D = "7.614";        E = "  +1.57%  " },
    @{ Row = 50; D = "9.877";        E = "  -1.50%  " },
    @{ Row = 51; D = "3.005";        E = "  +0.10%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
